# Tasks.xlsx update:
#  - Task 5 (row 6) and Task 6/ERD task (row 7) status moved from "Pending" to "Complete"
#  - Active selection moved to C9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the "Complete" formatting (bold green font, style index used by B2:B5)
# onto B6:B7, then set their text to "Complete" to match rows already marked done.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B6:B7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B6").Value2 = "Complete"
$ws.Range("B7").Value2 = "Complete"

# Move the active cell selection to C9, as last left by the author.
$ws.Range("C9").Select() | Out-Null
